$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# --- Players sheet: row-by-row corrections (player/team/game/status/stat line re-shuffle) ---
$ws.Range("D2").Value = 'Pablo Tamba'
$ws.Range("E2").Value = 'LSU'
$ws.Range("F2").Value = 'MSST@LSU'
$ws.Range("G2").Value = 'Final'
$ws.Range("H2").Value = 15
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 3
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 35
$ws.Range("D3").Value = 'Nate Ament'
$ws.Range("E3").Value = 'TENN'
$ws.Range("F3").Value = 'TENN@UGA'
$ws.Range("G3").Value = 'Final/OT'
$ws.Range("H3").Value = 13
$ws.Range("I3").Value = 19
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 2
$ws.Range("O3").Value = 38
$ws.Range("D4").Value = 'Matas Vokietaitis'
$ws.Range("H4").Value = 12
$ws.Range("I4").Value = 12
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 4
$ws.Range("O4").Value = 30
$ws.Range("C5").Value = 'Yes'
$ws.Range("D5").Value = 'Tramon Mark'
$ws.Range("E5").Value = 'TEX'
$ws.Range("F5").Value = 'TEX@AUB'
$ws.Range("H5").Value = -1
$ws.Range("I5").Value = 4
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 2
$ws.Range("O5").Value = 22
$ws.Range("C22").Value = 'Yes'
$ws.Range("D22").Value = 'Somtochukwu Cyril'
$ws.Range("E22").Value = 'UGA'
$ws.Range("F22").Value = 'TENN@UGA'
$ws.Range("G22").Value = 'Final/OT'
$ws.Range("H22").Value = 9
$ws.Range("I22").Value = 6
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2
$ws.Range("M22").Value = 3
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 32
$ws.Range("C23").Value = 'Yes'
$ws.Range("D23").Value = 'Mike Nwoko'
$ws.Range("E23").Value = 'LSU'
$ws.Range("F23").Value = 'MSST@LSU'
$ws.Range("G23").Value = 'Final'
$ws.Range("H23").Value = 3
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 17
$ws.Range("D24").Value = 'Kevin Overton'
$ws.Range("E24").Value = 'AUB'
$ws.Range("F24").Value = 'TEX@AUB'
$ws.Range("H24").Value = 25
$ws.Range("I24").Value = 25
$ws.Range("J24").Value = 1
$ws.Range("K24").Value = 1
$ws.Range("L24").Value = 1
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 33
$ws.Range("D25").Value = 'J.P. Estrella'
$ws.Range("E25").Value = 'TENN'
$ws.Range("H25").Value = 20
$ws.Range("I25").Value = 17
$ws.Range("J25").Value = 9
$ws.Range("K25").Value = 1
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("O25").Value = 31
$ws.Range("D26").Value = 'Xaivian Lee'
$ws.Range("E26").Value = 'FLA'
$ws.Range("F26").Value = 'FLA@SC'
$ws.Range("H26").Value = 14
$ws.Range("I26").Value = 6
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = 9
$ws.Range("L26").Value = 2
$ws.Range("O26").Value = 26
$ws.Range("D27").Value = 'Jordan Pope'
$ws.Range("E27").Value = 'TEX'
$ws.Range("F27").Value = 'TEX@AUB'
$ws.Range("H27").Value = 8
$ws.Range("I27").Value = 12
$ws.Range("J27").Value = 2
$ws.Range("K27").Value = 3
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = 2
$ws.Range("O27").Value = 35
$ws.Range("D30").Value = 'Marcus Millender'
$ws.Range("H30").Value = 22
$ws.Range("I30").Value = 19
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = 4
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = 1
$ws.Range("O30").Value = 28
$ws.Range("D31").Value = 'Kobe Knox'
$ws.Range("E31").Value = 'SC'
$ws.Range("F31").Value = 'FLA@SC'
$ws.Range("H31").Value = 7
$ws.Range("I31").Value = 6
$ws.Range("J31").Value = 4
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 1
$ws.Range("M31").Value = 1
$ws.Range("O31").Value = 26
$ws.Range("D32").Value = 'Blue Cain'
$ws.Range("H32").Value = 9
$ws.Range("I32").Value = 9
$ws.Range("J32").Value = 4
$ws.Range("K32").Value = 3
$ws.Range("L32").Value = 1
$ws.Range("N32").Value = 2
$ws.Range("O32").Value = 35
$ws.Range("D33").Value = 'Josh Hubbard'
$ws.Range("E33").Value = 'MSST'
$ws.Range("F33").Value = 'MSST@LSU'
$ws.Range("H33").Value = 9
$ws.Range("I33").Value = 15
$ws.Range("J33").Value = 2
$ws.Range("K33").Value = 3
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("O33").Value = 32

# --- OwnerTotals sheet: recomputed starter totals ---
$ws2.Range("B2").Value = 87
$ws2.Range("A4").Value = 'Booz'
$ws2.Range("B4").Value = 39
$ws2.Range("C4").Value = 4
$ws2.Range("A5").Value = 'Clay'
$ws2.Range("B5").Value = 38
$ws2.Range("C5").Value = 2
$ws2.Range("A6").Value = 'Hal'
$ws2.Range("B6").Value = 37
$ws2.Range("A7").Value = 'Mark'
$ws2.Range("B7").Value = 12
$ws2.Range("C7").Value = 2
$ws2.Range("A8").Value = 'Tar'
$ws2.Range("B8").Value = 3
$ws2.Range("C8").Value = 1
